$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD:AF) with the same formatting as the
# rest of the header row (bold, centered, thin box border) by copying the
# format from the last existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 79   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 83   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
